$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in boardings for existing rows (R5 Paoli-Thorndale PAO / 100 Norristown Speed Line)
$ws.Range("B23").Value = 576
$ws.Range("B24").Value = 294

# New Princeton Junction rail traffic section
$ws.Range("B41").Value = "Total"
$ws.Range("C41").Value = "Unique"

$ws.Range("A42").Value = "Princeton Junction"
$ws.Range("B42").Value = 6817

$ws.Hyperlinks.Add($ws.Range("D42"), "https://patch.com/new-jersey/livingston/here-are-new-jersey-transit-s-most-least-used-train-stations") | Out-Null
$ws.Range("D42").Style = "Hyperlink"

$ws.Range("B22").Select() | Out-Null
